$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dates in rows 31 and 32 (43924 -> 43925) ---
$ws.Range("C31").Value = 43925
$ws.Range("C32").Value = 43925

# --- Fill in previously-empty rows 33-35 with new activity-log entries ---
# Row 33
$ws.Range("B33").Value = 6977
$ws.Range("C33").Value = 43925
$ws.Range("D33").Value = 0.7090277777777777
$ws.Range("E33").Value = 0.75763888888888886
$ws.Range("G33").Value = "Added TestVectors. Discovered a bug in LogicUnit.vhd. Stopped for dinner"

# Row 34
$ws.Range("B34").Value = 6977
$ws.Range("C34").Value = 43925
$ws.Range("D34").Value = 0.88541666666666663
$ws.Range("E34").Value = 0.92708333333333337
$ws.Range("G34").Value = "Fixed LogicUnit.vhd."

# Row 35
$ws.Range("B35").Value = 6977
$ws.Range("C35").Value = 43925
$ws.Range("D35").Value = 0.93402777777777779
$ws.Range("E35").Value = 0.94791666666666663
$ws.Range("G35").Value = "Updated Functional Waveforms."

# Row 36 - only the "last 4 digits" column gets a value; rest stay blank
$ws.Range("B36").Value = 6977

# --- Update the view's scroll position / selection ---
$window = $excel.ActiveWindow
$window.ScrollRow = 22
$window.ScrollColumn = 1
$ws.Range("A36").Select()
